$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: insert the new "ownTeam"/"oppTeam" columns (D:E) and shift the
#     existing batsman..sr headers two columns to the right (F:K) ---
$ws.Range("D1").Value = "ownTeam"
$ws.Range("E1").Value = "oppTeam"
$ws.Range("F1").Value = "batsman"
$ws.Range("G1").Value = "totalRuns"
$ws.Range("H1").Value = "totalBalls"
$ws.Range("I1").Value = "total4s"
$ws.Range("J1").Value = "total6s"
$ws.Range("K1").Value = "sr"

# --- The totalRuns/totalBalls/total4s/total6s/sr columns hold figures such as
#     "0.00" / "175.00" / "-" that must stay literal text (not be reinterpreted
#     as numbers, which would drop the trailing zeros) ---
$ws.Range("G2:K8").NumberFormat = "@"

# Row 2 - v Kings XI Punjab, October 24 2020
$ws.Range("A2").Value = " Dubai (DSC)"
$ws.Range("B2").Value = " October 24 2020"
$ws.Range("C2").Value = "Kings XI won by 12 runs"
$ws.Range("D2").Value = "Sunrisers Hyderabad"
$ws.Range("E2").Value = "Kings XI Punjab"
$ws.Range("F2").Value = "Rashid Khan "
$ws.Range("G2").Value = "0"
$ws.Range("H2").Value = "1"
$ws.Range("I2").Value = "0"
$ws.Range("J2").Value = "0"
$ws.Range("K2").Value = "0.00"

# Row 3 - v Chennai Super Kings, October 13 2020
$ws.Range("A3").Value = " Dubai (DSC)"
$ws.Range("B3").Value = " October 13 2020"
$ws.Range("C3").Value = "Super Kings won by 20 runs"
$ws.Range("D3").Value = "Sunrisers Hyderabad"
$ws.Range("E3").Value = "Chennai Super Kings"
$ws.Range("F3").Value = "Rashid Khan "
$ws.Range("G3").Value = "14"
$ws.Range("H3").Value = "8"
$ws.Range("I3").Value = "1"
$ws.Range("J3").Value = "1"
$ws.Range("K3").Value = "175.00"

# Row 4 - v Kolkata Knight Riders, October 18 2020
$ws.Range("A4").Value = " Abu Dhabi"
$ws.Range("B4").Value = " October 18 2020"
$ws.Range("C4").Value = "Match tied (KKR won the one-over eliminator)"
$ws.Range("D4").Value = "Sunrisers Hyderabad"
$ws.Range("E4").Value = "Kolkata Knight Riders"
$ws.Range("F4").Value = "Rashid Khan "
$ws.Range("G4").Value = "1"
$ws.Range("H4").Value = "2"
$ws.Range("I4").Value = "0"
$ws.Range("J4").Value = "0"
$ws.Range("K4").Value = "50.00"

# Row 5 - v Royal Challengers Bangalore, September 21 2020
$ws.Range("A5").Value = " Dubai (DSC)"
$ws.Range("B5").Value = " September 21 2020"
$ws.Range("C5").Value = "RCB won by 10 runs"
$ws.Range("D5").Value = "Sunrisers Hyderabad"
$ws.Range("E5").Value = "Royal Challengers Bangalore"
$ws.Range("F5").Value = "Rashid Khan "
$ws.Range("G5").Value = "6"
$ws.Range("H5").Value = "5"
$ws.Range("I5").Value = "1"
$ws.Range("J5").Value = "0"
$ws.Range("K5").Value = "120.00"

# Row 6 - v Delhi Capitals, November 08 2020
$ws.Range("A6").Value = " Abu Dhabi"
$ws.Range("B6").Value = " November 08 2020"
$ws.Range("C6").Value = "Capitals won by 17 runs"
$ws.Range("D6").Value = "Sunrisers Hyderabad"
$ws.Range("E6").Value = "Delhi Capitals"
$ws.Range("F6").Value = "Rashid Khan "
$ws.Range("G6").Value = "11"
$ws.Range("H6").Value = "7"
$ws.Range("I6").Value = "1"
$ws.Range("J6").Value = "1"
$ws.Range("K6").Value = "157.14"

# Row 7 - v Kings XI Punjab, October 08 2020 (this was the only row already present)
$ws.Range("A7").Value = " Dubai (DSC)"
$ws.Range("B7").Value = " October 08 2020"
$ws.Range("C7").Value = "Sunrisers won by 69 runs"
$ws.Range("D7").Value = "Sunrisers Hyderabad"
$ws.Range("E7").Value = "Kings XI Punjab"
$ws.Range("F7").Value = "Rashid Khan "
$ws.Range("G7").Value = "0"
$ws.Range("H7").Value = "0"
$ws.Range("I7").Value = "0"
$ws.Range("J7").Value = "0"
$ws.Range("K7").Value = "-"

# Row 8 - v Mumbai Indians, October 04 2020
$ws.Range("A8").Value = " Sharjah"
$ws.Range("B8").Value = " October 04 2020"
$ws.Range("C8").Value = "Mumbai won by 34 runs"
$ws.Range("D8").Value = "Sunrisers Hyderabad"
$ws.Range("E8").Value = "Mumbai Indians"
$ws.Range("F8").Value = "Rashid Khan "
$ws.Range("G8").Value = "3"
$ws.Range("H8").Value = "7"
$ws.Range("I8").Value = "0"
$ws.Range("J8").Value = "0"
$ws.Range("K8").Value = "42.85"
